# Auto-generated edit script: updates market-data derived columns (H-N)
# per Sheets via scheduled runner (refreshed currentAveragePrice* and LevePrice*/LeveProfit* values).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$edits_ALC = @(
    @{Row=28; Col=8; Val=1487.8462},
    @{Row=28; Col=9; Val=1463.2222},
    @{Row=28; Col=11; Val=1463.2222},
    @{Row=28; Col=13; Val=-978.2221999999999},
    @{Row=33; Col=8; Val=4166786},
    @{Row=33; Col=9; Val=4545578},
    @{Row=33; Col=11; Val=4545578},
    @{Row=33; Col=13; Val=-4545349},
    @{Row=40; Col=8; Val=1189.1945},
    @{Row=40; Col=10; Val=1324.3334},
    @{Row=40; Col=12; Val=1324.3334},
    @{Row=40; Col=14; Val=-1674.3334},
    @{Row=43; Col=8; Val=14803.6},
    @{Row=43; Col=10; Val=9936.916999999999},
    @{Row=43; Col=12; Val=9936.916999999999},
    @{Row=43; Col=14; Val=-10074.917},
    @{Row=49; Col=8; Val=53.5},
    @{Row=49; Col=9; Val=17},
    @{Row=49; Col=10; Val=90},
    @{Row=49; Col=11; Val=51},
    @{Row=49; Col=12; Val=270},
    @{Row=49; Col=13; Val=85},
    @{Row=49; Col=14; Val=-542},
    @{Row=80; Col=8; Val=1602.5883},
    @{Row=80; Col=10; Val=1568.75},
    @{Row=80; Col=12; Val=4706.25},
    @{Row=80; Col=14; Val=-6702.25},
    @{Row=83; Col=8; Val=1602.5883},
    @{Row=83; Col=10; Val=1568.75},
    @{Row=83; Col=12; Val=14118.75},
    @{Row=83; Col=14; Val=-24102.75},
    @{Row=88; Col=8; Val=1823.5},
    @{Row=88; Col=9; Val=415.8},
    @{Row=88; Col=10; Val=2829},
    @{Row=88; Col=11; Val=415.8},
    @{Row=88; Col=12; Val=2829},
    @{Row=88; Col=13; Val=-9.800000000000011},
    @{Row=88; Col=14; Val=-3641},
    @{Row=91; Col=8; Val=1823.5},
    @{Row=91; Col=9; Val=415.8},
    @{Row=91; Col=10; Val=2829},
    @{Row=91; Col=11; Val=415.8},
    @{Row=91; Col=12; Val=2829},
    @{Row=91; Col=13; Val=988.2},
    @{Row=91; Col=14; Val=-5637}
)
foreach ($e in $edits_ALC) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("ARM")
$edits_ARM = @(
    @{Row=32; Col=8; Val=4294.1333},
    @{Row=32; Col=9; Val=3893.7144},
    @{Row=32; Col=11; Val=3893.7144},
    @{Row=32; Col=13; Val=-3606.7144},
    @{Row=61; Col=8; Val=2835},
    @{Row=61; Col=9; Val=2472.8262},
    @{Row=61; Col=11; Val=2472.8262},
    @{Row=61; Col=13; Val=-2260.8262},
    @{Row=74; Col=8; Val=38095.742},
    @{Row=74; Col=9; Val=44351.87},
    @{Row=74; Col=10; Val=2123},
    @{Row=74; Col=11; Val=44351.87},
    @{Row=74; Col=12; Val=2123},
    @{Row=74; Col=13; Val=-43477.87},
    @{Row=74; Col=14; Val=-3871},
    @{Row=77; Col=8; Val=38095.742},
    @{Row=77; Col=9; Val=44351.87},
    @{Row=77; Col=10; Val=2123},
    @{Row=77; Col=11; Val=221759.35},
    @{Row=77; Col=12; Val=10615},
    @{Row=77; Col=13; Val=-217391.35},
    @{Row=77; Col=14; Val=-19351},
    @{Row=88; Col=8; Val=2399.6897},
    @{Row=88; Col=9; Val=2301.6667},
    @{Row=88; Col=10; Val=2657},
    @{Row=88; Col=11; Val=2301.6667},
    @{Row=88; Col=12; Val=2657},
    @{Row=88; Col=13; Val=-1895.6667},
    @{Row=88; Col=14; Val=-3469},
    @{Row=91; Col=8; Val=2399.6897},
    @{Row=91; Col=9; Val=2301.6667},
    @{Row=91; Col=10; Val=2657},
    @{Row=91; Col=11; Val=2301.6667},
    @{Row=91; Col=12; Val=2657},
    @{Row=91; Col=13; Val=-897.6667000000002},
    @{Row=91; Col=14; Val=-5465},
    @{Row=125; Col=8; Val=32143},
    @{Row=125; Col=10; Val=32143},
    @{Row=125; Col=12; Val=32143},
    @{Row=125; Col=14; Val=-41983},
    @{Row=132; Col=8; Val=3604.0454},
    @{Row=132; Col=9; Val=3488.9473},
    @{Row=132; Col=11; Val=10466.8419},
    @{Row=132; Col=13; Val=-7936.841899999999},
    @{Row=136; Col=8; Val=2835},
    @{Row=136; Col=9; Val=2472.8262},
    @{Row=136; Col=11; Val=7418.4786},
    @{Row=136; Col=13; Val=-4868.4786}
)
foreach ($e in $edits_ARM) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("BSM")
$edits_BSM = @(
    @{Row=134; Col=8; Val=1903.7715},
    @{Row=134; Col=9; Val=1454.4},
    @{Row=134; Col=11; Val=4363.200000000001},
    @{Row=134; Col=13; Val=-1828.200000000001}
)
foreach ($e in $edits_BSM) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("CRP")
$edits_CRP = @(
    @{Row=31; Col=8; Val=89349.78999999999},
    @{Row=31; Col=9; Val=108836.4},
    @{Row=31; Col=10; Val=16275},
    @{Row=31; Col=11; Val=108836.4},
    @{Row=31; Col=12; Val=16275},
    @{Row=31; Col=13; Val=-108541.4},
    @{Row=31; Col=14; Val=-16865},
    @{Row=34; Col=8; Val=89349.78999999999},
    @{Row=34; Col=9; Val=108836.4},
    @{Row=34; Col=10; Val=16275},
    @{Row=34; Col=11; Val=108836.4},
    @{Row=34; Col=12; Val=16275},
    @{Row=34; Col=13; Val=-108634.4},
    @{Row=34; Col=14; Val=-16679},
    @{Row=50; Col=8; Val=0},
    @{Row=50; Col=10; Val=0},
    @{Row=50; Col=12; Val=0},
    @{Row=50; Col=14; Clear=$true},
    @{Row=60; Col=8; Val=4071.1428},
    @{Row=60; Col=9; Val=4071.1428},
    @{Row=60; Col=10; Val=0},
    @{Row=60; Col=11; Val=4071.1428},
    @{Row=60; Col=12; Val=0},
    @{Row=60; Col=14; Val=-3560.1428},
    @{Row=60; Col=13; Clear=$true},
    @{Row=62; Col=8; Val=5590.364},
    @{Row=62; Col=9; Val=6063},
    @{Row=62; Col=11; Val=6063},
    @{Row=62; Col=13; Val=-5439},
    @{Row=65; Col=8; Val=5590.364},
    @{Row=65; Col=9; Val=6063},
    @{Row=65; Col=11; Val=30315},
    @{Row=65; Col=13; Val=-27195},
    @{Row=122; Col=8; Val=1362.125},
    @{Row=122; Col=9; Val=1362.125},
    @{Row=122; Col=11; Val=4086.375},
    @{Row=122; Col=13; Val=-1636.375},
    @{Row=132; Col=8; Val=3800.7222},
    @{Row=132; Col=9; Val=3612.5293},
    @{Row=132; Col=11; Val=10837.5879},
    @{Row=132; Col=13; Val=-8307.5879}
)
foreach ($e in $edits_CRP) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("CUL")
$edits_CUL = @(
    @{Row=14; Col=8; Val=334488.66},
    @{Row=14; Col=9; Val=334488.66},
    @{Row=14; Col=11; Val=1003465.98},
    @{Row=14; Col=13; Val=-1003292.98},
    @{Row=33; Col=8; Val=220.11111},
    @{Row=33; Col=9; Val=48.333332},
    @{Row=33; Col=10; Val=306},
    @{Row=33; Col=11; Val=289.999992},
    @{Row=33; Col=12; Val=1836},
    @{Row=33; Col=13; Val=-6.99999200000002},
    @{Row=33; Col=14; Val=-2402},
    @{Row=69; Col=8; Val=20837},
    @{Row=69; Col=9; Val=1256},
    @{Row=69; Col=11; Val=3768},
    @{Row=69; Col=13; Val=-2957},
    @{Row=72; Col=8; Val=20837},
    @{Row=72; Col=9; Val=1256},
    @{Row=72; Col=11; Val=11304},
    @{Row=72; Col=13; Val=-7248},
    @{Row=107; Col=8; Val=25641544},
    @{Row=107; Col=9; Val=37037330},
    @{Row=107; Col=11; Val=111111990},
    @{Row=107; Col=13; Val=-111110070}
)
foreach ($e in $edits_CUL) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("GSM")
$edits_GSM = @(
    @{Row=80; Col=8; Val=2937.5},
    @{Row=80; Col=9; Val=2875},
    @{Row=80; Col=11; Val=2875},
    @{Row=80; Col=13; Val=-1877},
    @{Row=83; Col=8; Val=2937.5},
    @{Row=83; Col=9; Val=2875},
    @{Row=83; Col=11; Val=14375},
    @{Row=83; Col=13; Val=-9383},
    @{Row=92; Col=8; Val=13199.8},
    @{Row=92; Col=10; Val=13199.8},
    @{Row=92; Col=12; Val=13199.8},
    @{Row=92; Col=14; Val=-16943.8},
    @{Row=102; Col=8; Val=2139.182},
    @{Row=102; Col=9; Val=2139.182},
    @{Row=102; Col=11; Val=2139.182},
    @{Row=102; Col=13; Val=-517.1819999999998},
    @{Row=122; Col=8; Val=1450.7142},
    @{Row=122; Col=9; Val=1313},
    @{Row=122; Col=11; Val=3939},
    @{Row=122; Col=13; Val=-1489},
    @{Row=132; Col=8; Val=229756},
    @{Row=132; Col=9; Val=240506.28},
    @{Row=132; Col=11; Val=721518.84},
    @{Row=132; Col=13; Val=-718988.84}
)
foreach ($e in $edits_GSM) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("LTW")
$edits_LTW = @(
    @{Row=93; Col=8; Val=744440.25},
    @{Row=93; Col=9; Val=1012518.56},
    @{Row=93; Col=11; Val=1012518.56},
    @{Row=93; Col=13; Val=-1011270.56},
    @{Row=122; Col=8; Val=460176},
    @{Row=122; Col=9; Val=503291.16},
    @{Row=122; Col=10; Val=29024.5},
    @{Row=122; Col=11; Val=1509873.48},
    @{Row=122; Col=12; Val=87073.5},
    @{Row=122; Col=13; Val=-1507423.48},
    @{Row=122; Col=14; Val=-91973.5},
    @{Row=136; Col=8; Val=5144.6924},
    @{Row=136; Col=9; Val=4988.35},
    @{Row=136; Col=10; Val=5665.8335},
    @{Row=136; Col=11; Val=14965.05},
    @{Row=136; Col=12; Val=16997.5005},
    @{Row=136; Col=13; Val=-12415.05},
    @{Row=136; Col=14; Val=-22097.5005}
)
foreach ($e in $edits_LTW) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}

$ws = $wb.Worksheets.Item("WVR")
$edits_WVR = @(
    @{Row=49; Col=8; Val=0},
    @{Row=49; Col=10; Val=0},
    @{Row=49; Col=12; Val=0},
    @{Row=49; Col=14; Clear=$true},
    @{Row=70; Col=8; Val=37316.668},
    @{Row=70; Col=10; Val=37316.668},
    @{Row=70; Col=12; Val=37316.668},
    @{Row=70; Col=14; Val=-37946.668},
    @{Row=73; Col=8; Val=37316.668},
    @{Row=73; Col=10; Val=37316.668},
    @{Row=73; Col=12; Val=37316.668},
    @{Row=73; Col=14; Val=-39500.668},
    @{Row=81; Col=8; Val=7604.125},
    @{Row=81; Col=10; Val=3952.9412},
    @{Row=81; Col=12; Val=7905.8824},
    @{Row=81; Col=14; Val=-10027.8824},
    @{Row=84; Col=8; Val=7604.125},
    @{Row=84; Col=10; Val=3952.9412},
    @{Row=84; Col=12; Val=39529.412},
    @{Row=84; Col=14; Val=-50137.412},
    @{Row=136; Col=8; Val=1621.2},
    @{Row=136; Col=9; Val=1332.2307},
    @{Row=136; Col=11; Val=3996.6921},
    @{Row=136; Col=13; Val=-1446.6921},
    @{Row=138; Col=8; Val=70429},
    @{Row=138; Col=10; Val=70429},
    @{Row=138; Col=12; Val=70429},
    @{Row=138; Col=14; Val=-80709}
)
foreach ($e in $edits_WVR) {
    if ($e.ContainsKey("Clear")) {
        $ws.Cells.Item($e.Row, $e.Col).ClearContents()
    } else {
        $ws.Cells.Item($e.Row, $e.Col).Value = $e.Val
    }
}
